$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1945288753799392
$ws.Range("C2").Value = 0.5531914893617021
$ws.Range("J2").Value = 0.01519756838905775
$ws.Range("P2").Value = 0.1519756838905775
$ws.Range("S2").Value = 0.0851063829787234
$ws.Range("B3").Value = 0.005181347150259068
$ws.Range("C3").Value = 0.0155440414507772
$ws.Range("J3").Value = 0.02072538860103627
$ws.Range("P3").Value = 0.7564766839378239
$ws.Range("S3").Value = 0.2020725388601036
$ws.Range("J4").Value = 0.025
$ws.Range("P4").Value = 0.675
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.05019305019305019
$ws.Range("F6").Value = 0.05791505791505792
$ws.Range("J6").Value = 0.3011583011583012
$ws.Range("O6").Value = 0.01544401544401544
$ws.Range("Q6").Value = 0.1544401544401544
$ws.Range("R6").Value = 0.07722007722007722
$ws.Range("S6").Value = 0.3436293436293436
$ws.Range("B7").Value = 0.103448275862069
$ws.Range("D7").Value = 0.01149425287356322
$ws.Range("F7").Value = 0.07471264367816093
$ws.Range("J7").Value = 0.09195402298850575
$ws.Range("O7").Value = 0.01149425287356322
$ws.Range("Q7").Value = 0.1954022988505747
$ws.Range("R7").Value = 0.08620689655172414
$ws.Range("S7").Value = 0.4252873563218391
$ws.Range("B8").Value = 0.0872210953346856
$ws.Range("D8").Value = 0.02028397565922921
$ws.Range("E8").Value = 0.002028397565922921
$ws.Range("F8").Value = 0.07910750507099391
$ws.Range("J8").Value = 0.09533468559837728
$ws.Range("O8").Value = 0.02028397565922921
$ws.Range("Q8").Value = 0.18052738336714
$ws.Range("R8").Value = 0.08519269776876268
$ws.Range("S8").Value = 0.4300202839756592
$ws.Range("B9").Value = 0.0996309963099631
$ws.Range("D9").Value = 0.01476014760147601
$ws.Range("E9").Value = 0.003690036900369004
$ws.Range("F9").Value = 0.05166051660516605
$ws.Range("J9").Value = 0.08118081180811808
$ws.Range("O9").Value = 0.01476014760147601
$ws.Range("Q9").Value = 0.2140221402214022
$ws.Range("R9").Value = 0.0996309963099631
$ws.Range("S9").Value = 0.4206642066420664
$ws.Range("B10").Value = 0.1122956645344705
$ws.Range("D10").Value = 0.01705756929637527
$ws.Range("F10").Value = 0.06112295664534471
$ws.Range("J10").Value = 0.1066098081023454
$ws.Range("O10").Value = 0.01350390902629709
$ws.Range("Q10").Value = 0.2210376687988628
$ws.Range("R10").Value = 0.08528784648187633
$ws.Range("S10").Value = 0.3830845771144278
$ws.Range("G11").Value = 0.1235059760956175
$ws.Range("J11").Value = 0.09561752988047809
$ws.Range("K11").Value = 0.1912350597609562
$ws.Range("L11").Value = 0.5697211155378487
$ws.Range("S11").Value = 0.0199203187250996
$ws.Range("G12").Value = 0.7368421052631579
$ws.Range("J12").Value = 0.1842105263157895
$ws.Range("K12").Value = 0.01973684210526316
$ws.Range("L12").Value = 0.03947368421052631
$ws.Range("S12").Value = 0.01973684210526316
$ws.Range("G13").Value = 0.609375
$ws.Range("J13").Value = 0.359375
$ws.Range("S13").Value = 0.03125
$ws.Range("F15").Value = 0.01587301587301587
$ws.Range("H15").Value = 0.1825396825396825
$ws.Range("I15").Value = 0.08333333333333333
$ws.Range("J15").Value = 0.3214285714285715
$ws.Range("K15").Value = 0.03968253968253968
$ws.Range("M15").Value = 0.01587301587301587
$ws.Range("O15").Value = 0.07539682539682539
$ws.Range("S15").Value = 0.2658730158730159
$ws.Range("F16").Value = 0.02358490566037736
$ws.Range("H16").Value = 0.1415094339622641
$ws.Range("I16").Value = 0.09905660377358491
$ws.Range("J16").Value = 0.4669811320754717
$ws.Range("K16").Value = 0.07547169811320754
$ws.Range("M16").Value = 0.01886792452830189
$ws.Range("O16").Value = 0.05660377358490566
$ws.Range("S16").Value = 0.1179245283018868
$ws.Range("F17").Value = 0.02429906542056075
$ws.Range("H17").Value = 0.1700934579439252
$ws.Range("I17").Value = 0.102803738317757
$ws.Range("J17").Value = 0.4149532710280374
$ws.Range("K17").Value = 0.08224299065420561
$ws.Range("M17").Value = 0.02616822429906542
$ws.Range("N17").Value = 0.001869158878504673
$ws.Range("O17").Value = 0.06728971962616823
$ws.Range("S17").Value = 0.1102803738317757
$ws.Range("F18").Value = 0.02690582959641256
$ws.Range("H18").Value = 0.2017937219730942
$ws.Range("I18").Value = 0.1210762331838565
$ws.Range("J18").Value = 0.3901345291479821
$ws.Range("K18").Value = 0.06278026905829596
$ws.Range("M18").Value = 0.01345291479820628
$ws.Range("O18").Value = 0.06278026905829596
$ws.Range("S18").Value = 0.1210762331838565
$ws.Range("F19").Value = 0.02520435967302452
$ws.Range("H19").Value = 0.1900544959128065
$ws.Range("I19").Value = 0.1001362397820163
$ws.Range("J19").Value = 0.3712534059945504
$ws.Range("K19").Value = 0.08038147138964577
$ws.Range("M19").Value = 0.02861035422343324
$ws.Range("N19").Value = 0.0006811989100817438
$ws.Range("O19").Value = 0.07288828337874659
$ws.Range("S19").Value = 0.1307901907356948
